$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1966.6666
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H101").Value = 1294
$ws.Range("J101").Value = 1294
$ws.Range("L101").Value = 3882
$ws.Range("N101").Value = -7126
$ws.Range("H137").Value = 30906.646
$ws.Range("I137").Value = 1456
$ws.Range("K137").Value = 4368
$ws.Range("M137").Value = -1818
$ws.Range("H138").Value = 2170.5205
$ws.Range("J138").Value = 2300.7817
$ws.Range("L138").Value = 6902.3451
$ws.Range("N138").Value = -17182.3451

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 33698.76
$ws.Range("I32").Value = 39147.188
$ws.Range("K32").Value = 39147.188
$ws.Range("M32").Value = -38860.188
$ws.Range("H74").Value = 40001310
$ws.Range("I74").Value = 55556308
$ws.Range("J74").Value = 2750.7144
$ws.Range("K74").Value = 55556308
$ws.Range("L74").Value = 2750.7144
$ws.Range("M74").Value = -55555434
$ws.Range("N74").Value = -4498.7144
$ws.Range("H77").Value = 40001310
$ws.Range("I77").Value = 55556308
$ws.Range("J77").Value = 2750.7144
$ws.Range("K77").Value = 277781540
$ws.Range("L77").Value = 13753.572
$ws.Range("M77").Value = -277777172
$ws.Range("N77").Value = -22489.572
$ws.Range("H97").Value = 867.4286
$ws.Range("I97").Value = 940.3
$ws.Range("J97").Value = 685.25
$ws.Range("K97").Value = 940.3
$ws.Range("L97").Value = 685.25
$ws.Range("M97").Value = -444.3
$ws.Range("N97").Value = -1677.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 8691.333000000001
$ws.Range("I36").Value = 3037
$ws.Range("J36").Value = 20000
$ws.Range("K36").Value = 3037
$ws.Range("L36").Value = 20000
$ws.Range("M36").Value = -2503
$ws.Range("N36").Value = -21068
$ws.Range("H58").Value = 19800
$ws.Range("J58").Value = 19800
$ws.Range("L58").Value = 19800
$ws.Range("N58").Value = -20388
$ws.Range("H60").Value = 19970
$ws.Range("J60").Value = 19970
$ws.Range("L60").Value = 19970
$ws.Range("N60").Value = -21168
$ws.Range("H139").Value = 48880
$ws.Range("J139").Value = 48880
$ws.Range("L139").Value = 48880
$ws.Range("N139").Value = -59160

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 70.666664
$ws.Range("I7").Value = 30
$ws.Range("J7").Value = 91
$ws.Range("K7").Value = 30
$ws.Range("L7").Value = 91
$ws.Range("M7").Value = 83
$ws.Range("N7").Value = -317
$ws.Range("H120").Value = 9800
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H132").Value = 17038.363
$ws.Range("I132").Value = 18938.55
$ws.Range("K132").Value = 56815.64999999999
$ws.Range("M132").Value = -54285.64999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 3548.2
$ws.Range("I64").Value = 1653
$ws.Range("K64").Value = 4959
$ws.Range("M64").Value = -4689
$ws.Range("H67").Value = 3548.2
$ws.Range("I67").Value = 1653
$ws.Range("K67").Value = 4959
$ws.Range("M67").Value = -4023
$ws.Range("H122").Value = 674.86957
$ws.Range("I122").Value = 324.5
$ws.Range("J122").Value = 944.38464
$ws.Range("K122").Value = 2920.5
$ws.Range("L122").Value = 8499.46176
$ws.Range("M122").Value = -470.5
$ws.Range("N122").Value = -13399.46176
$ws.Range("H131").Value = 527152.0600000001
$ws.Range("J131").Value = 527152.0600000001
$ws.Range("L131").Value = 1581456.18
$ws.Range("N131").Value = -1591536.18

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 79.22221999999999
$ws.Range("I2").Value = 69.666664
$ws.Range("J2").Value = 98.333336
$ws.Range("K2").Value = 69.666664
$ws.Range("L2").Value = 98.333336
$ws.Range("M2").Value = 43.333336
$ws.Range("N2").Value = -324.333336
$ws.Range("H18").Value = 10005
$ws.Range("I18").Value = 10005
$ws.Range("K18").Value = 10005
$ws.Range("M18").Value = -9712
$ws.Range("H44").Value = 20000
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 20000
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 20000
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -21192
$ws.Range("H46").Value = 25600
$ws.Range("J46").Value = 25600
$ws.Range("L46").Value = 25600
$ws.Range("N46").Value = -25912
$ws.Range("H52").Value = 20005600
$ws.Range("J52").Value = 20005600
$ws.Range("L52").Value = 20005600
$ws.Range("N52").Value = -20006118
$ws.Range("H57").Value = 28695
$ws.Range("J57").Value = 29990
$ws.Range("L57").Value = 29990
$ws.Range("N57").Value = -31630
$ws.Range("H70").Value = 11359.934
$ws.Range("I70").Value = 19866.666
$ws.Range("J70").Value = 5688.778
$ws.Range("K70").Value = 19866.666
$ws.Range("L70").Value = 5688.778
$ws.Range("M70").Value = -19596.666
$ws.Range("N70").Value = -6228.778
$ws.Range("H73").Value = 11359.934
$ws.Range("I73").Value = 19866.666
$ws.Range("J73").Value = 5688.778
$ws.Range("K73").Value = 19866.666
$ws.Range("L73").Value = 5688.778
$ws.Range("M73").Value = -18930.666
$ws.Range("N73").Value = -7560.778
$ws.Range("H80").Value = 3963.4167
$ws.Range("I80").Value = 3271
$ws.Range("J80").Value = 4458
$ws.Range("K80").Value = 3271
$ws.Range("L80").Value = 4458
$ws.Range("M80").Value = -2273
$ws.Range("N80").Value = -6454
$ws.Range("H83").Value = 3963.4167
$ws.Range("I83").Value = 3271
$ws.Range("J83").Value = 4458
$ws.Range("K83").Value = 16355
$ws.Range("L83").Value = 22290
$ws.Range("M83").Value = -11363
$ws.Range("N83").Value = -32274
$ws.Range("H107").Value = 9615738
$ws.Range("I107").Value = 421.66666
$ws.Range("K107").Value = 421.66666
$ws.Range("M107").Value = 1498.33334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H22").Value = 2826.35
$ws.Range("I22").Value = 3556.6667
$ws.Range("J22").Value = 2228.818
$ws.Range("K22").Value = 3556.6667
$ws.Range("L22").Value = 2228.818
$ws.Range("M22").Value = -3261.6667
$ws.Range("N22").Value = -2818.818
$ws.Range("H27").Value = 2826.35
$ws.Range("I27").Value = 3556.6667
$ws.Range("J27").Value = 2228.818
$ws.Range("K27").Value = 3556.6667
$ws.Range("L27").Value = 2228.818
$ws.Range("M27").Value = -3449.6667
$ws.Range("N27").Value = -2442.818
$ws.Range("H41").Value = 15000
$ws.Range("J41").Value = 15000
$ws.Range("L41").Value = 15000
$ws.Range("N41").Value = -15876
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("N48").ClearContents()
$ws.Range("H132").Value = 1643.125
$ws.Range("I132").Value = 1123.28
$ws.Range("J132").Value = 3499.7144
$ws.Range("K132").Value = 3369.84
$ws.Range("L132").Value = 10499.1432
$ws.Range("M132").Value = -839.8400000000001
$ws.Range("N132").Value = -15559.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 3333.1667
$ws.Range("J15").Value = 3333.1667
$ws.Range("L15").Value = 3333.1667
$ws.Range("N15").Value = -3909.1667
$ws.Range("H21").Value = 1965
$ws.Range("J21").Value = 1965
$ws.Range("L21").Value = 1965
$ws.Range("N21").Value = -2435
$ws.Range("H35").Value = 1965
$ws.Range("J35").Value = 1965
$ws.Range("L35").Value = 1965
$ws.Range("N35").Value = -2545
$ws.Range("H37").Value = 41264.5
$ws.Range("I37").Value = 10000
$ws.Range("J37").Value = 51686
$ws.Range("K37").Value = 10000
$ws.Range("L37").Value = 51686
$ws.Range("M37").Value = -9797
$ws.Range("N37").Value = -52092
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H107").Value = 5682793
$ws.Range("I107").Value = 1350
$ws.Range("J107").Value = 11364236
$ws.Range("K107").Value = 4050
$ws.Range("L107").Value = 34092708
$ws.Range("M107").Value = -2130
$ws.Range("N107").Value = -34096548
